# Add DBWriter Functionality: Get Batch Number
#
# The sheet used to start with a header row (Cno | Course Code | Course
# Title | Batch | Teacher Name). That header row - and the highlighted
# fill/border formatting that only it used - is no longer wanted, so we
# delete row 1 outright. This shifts every data row up by one (the former
# row 2 becomes row 1, etc.) and leaves the plain, unformatted data table
# behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Rows("1").Delete()

# Match the saved view state: selection on C6, scrolled back to the top.
$ws.Range("C6").Select()
